# Automatische test-sync: 2025-06-19 08:00:10
# Append a new log entry to the "Logs" sheet and update the corresponding
# category count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 19

$logs.Range("A$newRow").Value = "Vragen over samenwerking"
$logs.Range("B$newRow").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$newRow").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D$newRow").Value = "Overig"
$logs.Range("F$newRow").Value = "2025-06-18 16:30:10"
$logs.Range("G$newRow").Value = "Nee"

# Update the "Overig" tally on the Dashboard sheet.
$dashboard.Range("B2").Value = 9

# Extend the conditional formatting ranges so the new row is covered too.
foreach ($fc in $logs.Range("D2").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("D2:D$newRow"))
}
foreach ($fc in $logs.Range("G2").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("G2:G$newRow"))
}
